$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.967.40"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "1.859.31"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.84"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4359"
$ws.Range("E7").Value = "  -4.90%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07494"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9409"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.30"
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("D12").Value = "1.852.30"
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.723"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.436"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06855"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.58"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009009"
$ws.Range("E18").Value = "  -4.75%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.97"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("D21").Value = "27.937.62"
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "2.080.29"
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.007"
$ws.Range("E25").Value = "  -4.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.56"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.386"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.80"
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.733"
$ws.Range("E30").Value = "  -7.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8064"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.827"
$ws.Range("E33").Value = "  -4.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.175"
$ws.Range("E34").Value = "  -5.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.949"
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05483"
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("E38").Value = "  -3.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01978"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.921"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5257"
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.025"
$ws.Range("E42").Value = "  -5.24%  "
$ws.Range("E43").Value = "  -4.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.799"
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06816"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4904"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.59"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.22"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.918"
$ws.Range("E49").Value = "  -9.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.679"
$ws.Range("E50").Value = "  -5.46%  "
$ws.Range("E51").Value = "  -0.18%  "
